# project started day1 first commit
# Add a new data row (Hello / Hellodesc) below the existing header row,
# then move the active selection to C7 as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Hello"
$ws.Range("B2").Value = "Hellodesc"

[void]$ws.Range("C7").Select()
